$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename shared-string cluster labels (sending cluster col A, target cluster col D)
# "Neutrophils" -> "MuSCs" and former "MuSCs" -> "Inflammatory-Mac"

# Row 2
$ws.Range("G2").Value = 0.1806205
$ws.Range("H2").Value = 0.361241
$ws.Range("I2").Value = 0.7284509268949775
$ws.Range("J2").Value = 0.7284509268949775
$ws.Range("M2").Value = 6.9192795
$ws.Range("N2").Value = 13.838559
$ws.Range("O2").Value = 0.6281451873364243
$ws.Range("P2").Value = 0.5718408792853329
$ws.Range("Q2").Value = 1.24976372292975
$ws.Range("R2").Value = 4.999054891719
$ws.Range("S2").Value = 0.4575729439398376
$ws.Range("T2").Value = 0.4165580185518397

# Row 3
$ws.Range("G3").Value = 0.1806205
$ws.Range("H3").Value = 0.361241
$ws.Range("I3").Value = 0.7284509268949775
$ws.Range("J3").Value = 0.7284509268949775
$ws.Range("O3").Value = 0.1890685122662809
$ws.Range("P3").Value = 0.2581818021036928
$ws.Range("Q3").Value = 0.376172535494
$ws.Range("R3").Value = 2.257035212964
$ws.Range("S3").Value = 0.1377271330070268
$ws.Range("T3").Value = 0.1880727730498507

# Row 4
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("G4").Value = 0.1806205
$ws.Range("H4").Value = 0.361241
$ws.Range("I4").Value = 0.7284509268949775
$ws.Range("J4").Value = 0.7284509268949775
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.07724199999999999
$ws.Range("N4").Value = 0.231726
$ws.Range("O4").Value = 0.007012173819577614
$ws.Range("P4").Value = 0.009575447818900294
$ws.Range("Q4").Value = 0.013951488661
$ws.Range("R4").Value = 0.08370893196599999
$ws.Range("S4").Value = 0.005108024518420008
$ws.Range("T4").Value = 0.00697524383911241

# Row 5
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 0.1806205
$ws.Range("H5").Value = 0.361241
$ws.Range("I5").Value = 0.7284509268949775
$ws.Range("J5").Value = 0.7284509268949775
$ws.Range("K5").Value = 2
$ws.Range("M5").Value = 1.9269465
$ws.Range("N5").Value = 3.853893
$ws.Range("O5").Value = 0.1749318220531151
$ws.Range("P5").Value = 0.1592516649884999
$ws.Range("Q5").Value = 0.34804604030325
$ws.Range("R5").Value = 1.392184161213
$ws.Range("S5").Value = 0.127429247918019
$ws.Range("T5").Value = 0.1160070229704411

# Row 6
$ws.Range("G6").Value = 0.1806205
$ws.Range("H6").Value = 0.361241
$ws.Range("I6").Value = 0.7284509268949775
$ws.Range("J6").Value = 0.7284509268949775
$ws.Range("M6").Value = 0.009278333333333333
$ws.Range("N6").Value = 0.027835
$ws.Range("O6").Value = 0.0008423045246020856
$ws.Range("P6").Value = 0.001150205803574436
$ws.Range("Q6").Value = 0.001675857205833333
$ws.Range("R6").Value = 0.010055143235
$ws.Range("S6").Value = 0.0006135775116742227
$ws.Range("T6").Value = 0.00083786848373378

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.067331
$ws.Range("H7").Value = 0.134662
$ws.Range("I7").Value = 0.2715490731050226
$ws.Range("J7").Value = 0.2715490731050226
$ws.Range("M7").Value = 6.9192795
$ws.Range("N7").Value = 13.838559
$ws.Range("O7").Value = 0.6281451873364243
$ws.Range("P7").Value = 0.5718408792853329
$ws.Range("Q7").Value = 0.4658820080145
$ws.Range("R7").Value = 1.863528032058
$ws.Range("S7").Value = 0.1705722433965868
$ws.Range("T7").Value = 0.1552828607334933

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.067331
$ws.Range("H8").Value = 0.134662
$ws.Range("I8").Value = 0.2715490731050226
$ws.Range("J8").Value = 0.2715490731050226
$ws.Range("O8").Value = 0.1890685122662809
$ws.Range("P8").Value = 0.2581818021036928
$ws.Range("Q8").Value = 0.140228119108
$ws.Range("R8").Value = 0.841368714648
$ws.Range("S8").Value = 0.05134137925925418
$ws.Range("T8").Value = 0.07010902905384214

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = 0.067331
$ws.Range("H9").Value = 0.134662
$ws.Range("I9").Value = 0.2715490731050226
$ws.Range("J9").Value = 0.2715490731050226
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.07724199999999999
$ws.Range("N9").Value = 0.231726
$ws.Range("O9").Value = 0.007012173819577614
$ws.Range("P9").Value = 0.009575447818900294
$ws.Range("Q9").Value = 0.005200781101999999
$ws.Range("R9").Value = 0.031204686612
$ws.Range("S9").Value = 0.001904149301157607
$ws.Range("T9").Value = 0.002600203979787885

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("F10").Value = 0.5
$ws.Range("G10").Value = 0.067331
$ws.Range("H10").Value = 0.134662
$ws.Range("I10").Value = 0.2715490731050226
$ws.Range("J10").Value = 0.2715490731050226
$ws.Range("K10").Value = 2
$ws.Range("M10").Value = 1.9269465
$ws.Range("N10").Value = 3.853893
$ws.Range("O10").Value = 0.1749318220531151
$ws.Range("P10").Value = 0.1592516649884999
$ws.Range("Q10").Value = 0.1297432347915
$ws.Range("R10").Value = 0.5189729391660001
$ws.Range("S10").Value = 0.04750257413509616
$ws.Range("T10").Value = 0.04324464201805871

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("F11").Value = 0.5
$ws.Range("G11").Value = 0.067331
$ws.Range("H11").Value = 0.134662
$ws.Range("I11").Value = 0.2715490731050226
$ws.Range("J11").Value = 0.2715490731050226
$ws.Range("M11").Value = 0.009278333333333333
$ws.Range("N11").Value = 0.027835
$ws.Range("O11").Value = 0.0008423045246020856
$ws.Range("P11").Value = 0.001150205803574436
$ws.Range("Q11").Value = 0.0006247194616666666
$ws.Range("R11").Value = 0.00374831677
$ws.Range("S11").Value = 0.000228727012927863
$ws.Range("T11").Value = 0.0003123373198406556

